$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "68.826.27"
$ws.Cells.Item(2, 5).Value = "  +0.41%  "

$ws.Cells.Item(3, 4).Value = "2.462.77"
$ws.Cells.Item(3, 5).Value = "  +0.68%  "

$ws.Cells.Item(4, 5).Value = "  +0.00%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "558.95"
$ws.Cells.Item(5, 5).Value = "  -0.71%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "161.68"
$ws.Cells.Item(6, 5).Value = "  -0.49%  "

$ws.Cells.Item(7, 5).Value = "  +0.00%  "

$ws.Cells.Item(8, 5).Value = "  -0.05%  "

$ws.Cells.Item(9, 4).Value = "2.461.63"
$ws.Cells.Item(9, 5).Value = "  +0.65%  "

$ws.Cells.Item(10, 5).Value = "  -0.71%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.165"
$ws.Cells.Item(11, 5).Value = "  +0.58%  "

$ws.Cells.Item(12, 5).Value = "  +1.00%  "

$ws.Cells.Item(13, 5).Value = "  -3.48%  "

$ws.Cells.Item(14, 5).Value = "  +0.26%  "

$ws.Cells.Item(15, 4).Value = "68.751.58"
$ws.Cells.Item(15, 5).Value = "  +0.52%  "

$ws.Cells.Item(16, 5).Value = "  -1.75%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "23.54"
$ws.Cells.Item(17, 5).Value = "  -0.47%  "

$ws.Cells.Item(18, 4).Value = "2.454.34"
$ws.Cells.Item(18, 5).Value = "  +0.38%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "10.68"
$ws.Cells.Item(19, 5).Value = "  -2.65%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "335.10"
$ws.Cells.Item(20, 5).Value = "  -3.03%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.93"
$ws.Cells.Item(21, 5).Value = "  -3.43%  "

$ws.Cells.Item(22, 5).Value = "  -0.83%  "

$ws.Cells.Item(23, 2).Value = "Dai"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "1.00"
$ws.Cells.Item(23, 5).Value = "  +0.05%  "

$ws.Cells.Item(24, 2).Value = "SuiNetwork"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "1.88"
$ws.Cells.Item(24, 5).Value = "  +0.26%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "66.84"
$ws.Cells.Item(25, 5).Value = "  -2.10%  "

$ws.Cells.Item(26, 4).Value = "2.600.69"

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "3.65"
$ws.Cells.Item(27, 5).Value = "  -2.79%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.00"
$ws.Cells.Item(28, 5).Value = "  -2.72%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "8.15"
$ws.Cells.Item(29, 5).Value = "  -1.14%  "

$ws.Cells.Item(30, 4).Value = "0.0₃0813"
$ws.Cells.Item(30, 5).Value = "  -2.97%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "7.16"
$ws.Cells.Item(31, 5).Value = "  -1.97%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.999"
$ws.Cells.Item(32, 5).Value = "  -0.06%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "427.72"
$ws.Cells.Item(33, 5).Value = "  -0.47%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.13"
$ws.Cells.Item(34, 5).Value = "  -3.38%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.61"
$ws.Cells.Item(35, 5).Value = "  -3.71%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "157.96"
$ws.Cells.Item(36, 5).Value = "  +0.95%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "19.03"
$ws.Cells.Item(37, 5).Value = "  +0.15%  "

$ws.Cells.Item(38, 5).Value = "  -0.05%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.108"
$ws.Cells.Item(39, 5).Value = "  -1.27%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "17.74"
$ws.Cells.Item(40, 5).Value = "  -1.12%  "

$ws.Cells.Item(41, 5).Value = "  -2.30%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "4.40"

$ws.Cells.Item(43, 5).Value = "  -4.77%  "

$ws.Cells.Item(44, 5).Value = "  -2.27%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.04"
$ws.Cells.Item(45, 5).Value = "  -1.91%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "131.87"
$ws.Cells.Item(46, 5).Value = "  -2.34%  "

$ws.Cells.Item(47, 5).Value = "  -0.98%  "

$ws.Cells.Item(48, 5).Value = "  -0.58%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.482"
$ws.Cells.Item(49, 5).Value = "  -1.69%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.559"
$ws.Cells.Item(50, 5).Value = "  -0.44%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0910"
$ws.Cells.Item(51, 5).Value = "  -0.46%  "
